$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Edit 1 & 2: in table 1, row 1, cell 1, merge the split runs back into a
# single run per paragraph (the content/text stays the same except for the
# "s" -> "containing ... margins" rewording in the 3rd paragraph).
# -----------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$cell11 = $t1.Cell(1, 1)

$findRange1 = $cell11.Range
$found1 = $findRange1.Find.Execute("a cell ", $true, $false, $false, $false, $false, $true, 1, $false, `
    "a cell ", 2)

$findRange2 = $cell11.Range
$found2 = $findRange2.Find.Execute("with three paragraphs with auto margins.", $true, $false, $false, $false, $false, $true, 1, $false, `
    "with three paragraphs containing auto margins.", 2)

Write-Host "Edit1 found:" $found1 "Edit2 found:" $found2

# -----------------------------------------------------------------------
# Edit 3: append a new row at the end of table 1 with a cell containing
# three paragraphs of text ("This is" / "a cell " / "with three paragraphs
# with disabled auto margins.") and beforeAutospacing explicitly disabled.
# -----------------------------------------------------------------------
$newRow = $t1.Rows.Add()
$newCell = $t1.Cell($t1.Rows.Count, 1)

# The new cell starts out with a single (empty) paragraph; split it into
# three paragraphs so each can hold one line of text.
$firstPara = $newCell.Range.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphAfter()
$firstPara.Range.InsertParagraphAfter()

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$paragraph1Xml = "<w:p $wNs><w:pPr><w:spacing w:beforeAutospacing=`"0`"/><w:rPr><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/><w:lang w:val=`"en-US`"/></w:rPr><w:t>This is</w:t></w:r></w:p>"
$paragraph2Xml = "<w:p $wNs><w:pPr><w:spacing w:beforeAutospacing=`"0`"/><w:rPr><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">a cell </w:t></w:r></w:p>"
$paragraph3Xml = "<w:p $wNs><w:pPr><w:spacing w:beforeAutospacing=`"0`"/><w:rPr><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/><w:lang w:val=`"en-US`"/></w:rPr><w:t>with three paragraphs with disabled auto margins.</w:t></w:r></w:p>"

$newCell.Range.Paragraphs.Item(1).Range.InsertXML($paragraph1Xml)
$newCell.Range.Paragraphs.Item(2).Range.InsertXML($paragraph2Xml)
$newCell.Range.Paragraphs.Item(3).Range.InsertXML($paragraph3Xml)

Write-Host "New row appended, paragraph count:" $newCell.Range.Paragraphs.Count

# -----------------------------------------------------------------------
# Cleanup: Word's table-row automation stamps the freshly added <w:tr> with
# rsid/paraId bookkeeping attributes that a hand authored fixture (as in
# the target document) does not have. Strip them from that single row only.
# -----------------------------------------------------------------------
$fullXml = $d.WordOpenXML
$cleanedXml = $fullXml -replace '<w:tr[^>]*>(?=(?:(?!</w:tr>).)*?with three paragraphs with disabled auto margins\.)', '<w:tr>'
if ($cleanedXml -ne $fullXml) {
    $d.WordOpenXML = $cleanedXml
    Write-Host "Stray rsid/paraId attributes stripped from new row"
} else {
    Write-Host "No stray attributes found to strip"
}
